$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    if ($bCell.Value2 -eq "bleu") {
        $bCell.Value = "noir"
    }

    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq "résultat et / ou publication posté") {
        $cCell.Value = "résultat postés ou publiés"
    }
    elseif ($cCell.Value2 -eq "pas de résultat ni de publication") {
        $cCell.Value = "pas de résultat postés ni publiés"
    }
    elseif ($cCell.Value2 -eq "résultat et / ou publication posté dans les 36 mois") {
        $cCell.Value = "résultat postés ou publiés dans les 36 mois"
    }
}
